$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 254
$ws.Range("I38").Value = 254
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 762
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -390
$ws.Range("H58").Value = 361.5
$ws.Range("I58").Value = 233.8
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 701.4000000000001
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -551.4000000000001
$ws.Range("N58").Value = -3300
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H132").Value = 1899.75
$ws.Range("I132").Value = 1899.6666
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 5698.9998
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -3168.9998
$ws.Range("N132").Value = -10760
$ws.Range("H137").Value = 4278.3
$ws.Range("I137").Value = 1917.6
$ws.Range("J137").Value = 6639
$ws.Range("K137").Value = 5752.799999999999
$ws.Range("L137").Value = 19917
$ws.Range("M137").Value = -3202.799999999999
$ws.Range("N137").Value = -25017

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2149.5715
$ws.Range("I74").Value = 2209.6
$ws.Range("J74").Value = 1999.5
$ws.Range("K74").Value = 2209.6
$ws.Range("L74").Value = 1999.5
$ws.Range("M74").Value = -1335.6
$ws.Range("N74").Value = -3747.5
$ws.Range("H77").Value = 2149.5715
$ws.Range("I77").Value = 2209.6
$ws.Range("J77").Value = 1999.5
$ws.Range("K77").Value = 11048
$ws.Range("L77").Value = 9997.5
$ws.Range("M77").Value = -6680
$ws.Range("N77").Value = -18733.5
$ws.Range("H132").Value = 2928.4
$ws.Range("I132").Value = 2928.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8785.200000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6255.200000000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1237
$ws.Range("I105").Value = 1237
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1237
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 510
$ws.Range("H134").Value = 6351.2
$ws.Range("I134").Value = 6351.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 19053.6
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -16518.6

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3166
$ws.Range("I31").Value = 2192.4
$ws.Range("J31").Value = 5600
$ws.Range("K31").Value = 2192.4
$ws.Range("L31").Value = 5600
$ws.Range("M31").Value = -1897.4
$ws.Range("N31").Value = -6190
$ws.Range("H34").Value = 3166
$ws.Range("I34").Value = 2192.4
$ws.Range("J34").Value = 5600
$ws.Range("K34").Value = 2192.4
$ws.Range("L34").Value = 5600
$ws.Range("M34").Value = -1990.4
$ws.Range("N34").Value = -6004
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10368
$ws.Range("H101").Value = 10000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 10000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 10000
$ws.Range("N101").Value = -16490
$ws.Range("H111").Value = 70000
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 70000
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 70000
$ws.Range("N111").Value = -78180
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("G129").Value = 35378
$ws.Range("H129").Value = 58000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 58000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 58000
$ws.Range("N129").Value = -68000
$ws.Range("G130").Value = 34689
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("G131").Value = 35461
$ws.Range("H131").Value = 50000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 50000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("G132").Value = 44019
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060
$ws.Range("G133").Value = 43328
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("G134").Value = 44020
$ws.Range("H134").Value = 961.38464
$ws.Range("I134").Value = 927.0909
$ws.Range("J134").Value = 1150
$ws.Range("K134").Value = 2781.2727
$ws.Range("L134").Value = 3450
$ws.Range("M134").Value = -246.2727
$ws.Range("N134").Value = -8520
$ws.Range("G135").Value = 42008
$ws.Range("H135").Value = 68618.5
$ws.Range("I135").Value = 49666.668
$ws.Range("J135").Value = 79989.60000000001
$ws.Range("K135").Value = 49666.668
$ws.Range("L135").Value = 79989.60000000001
$ws.Range("M135").Value = -44596.668
$ws.Range("N135").Value = -90129.60000000001
$ws.Range("G137").Value = 43231
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("G138").Value = 42302
$ws.Range("H138").Value = 60000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 60000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("G140").Value = 42455
$ws.Range("H140").Value = 150000
$ws.Range("I140").Value = 150000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 150000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -144820
$ws.Range("G141").Value = 43345
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 285810.28
$ws.Range("I11").Value = 666733.3
$ws.Range("J11").Value = 118
$ws.Range("K11").Value = 2000199.9
$ws.Range("L11").Value = 354
$ws.Range("M11").Value = -2000059.9
$ws.Range("N11").Value = -634
$ws.Range("H16").Value = 447.25
$ws.Range("I16").Value = 429.66666
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1288.99998
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -1115.99998
$ws.Range("N16").Value = -1846
$ws.Range("H26").Value = 88.8
$ws.Range("I26").Value = 75
$ws.Range("J26").Value = 109.5
$ws.Range("K26").Value = 225
$ws.Range("L26").Value = 328.5
$ws.Range("M26").Value = 63
$ws.Range("N26").Value = -904.5
$ws.Range("H93").Value = 327
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 327
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 981
$ws.Range("N93").Value = -4725
$ws.Range("H115").Value = 2025.6666
$ws.Range("I115").Value = 999
$ws.Range("J115").Value = 2231
$ws.Range("K115").Value = 2997
$ws.Range("L115").Value = 6693
$ws.Range("M115").Value = -1822
$ws.Range("N115").Value = -9043
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1503.9
$ws.Range("I113").Value = 1120
$ws.Range("J113").Value = 2399.6667
$ws.Range("K113").Value = 1120
$ws.Range("L113").Value = 2399.6667
$ws.Range("M113").Value = 1050
$ws.Range("N113").Value = -6739.6667
$ws.Range("H132").Value = 1499
$ws.Range("I132").Value = 1499
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4497
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1967

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1532.3334
$ws.Range("I82").Value = 1348.875
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 1348.875
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -987.875
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 1532.3334
$ws.Range("I85").Value = 1348.875
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1348.875
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -100.875
$ws.Range("N85").Value = -5496
$ws.Range("H127").Value = 33333
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 33333
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 33333
$ws.Range("N127").Value = -43253

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1432
$ws.Range("I107").Value = 821.7143
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 2465.1429
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -545.1428999999998
$ws.Range("N107").Value = -11340
$ws.Range("H132").Value = 2125.75
$ws.Range("I132").Value = 1715.1428
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5145.428400000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2615.428400000001
$ws.Range("N132").Value = -20060
